$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.08550000000001
$ws.Range("B4").Value = 4.942300000000007

$ws.Range("B5").Value = 5.273399999999999

$ws.Range("A6").Value = -21.35290000000001

$ws.Range("A7").Value = -21.531

$ws.Range("B8").Value = 5.1273

$ws.Range("A16").Value = -21.60510000000001
$ws.Range("B16").Value = 4.810500000000002

$ws.Range("A20").Value = -22.12060000000003

$ws.Range("B22").Value = 5.189900000000005
